$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new year column S (2022), copying the formatting from the adjacent
# existing year columns (R4 header style, R5 data style) before writing values.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 13.5

# Update existing 2019/2020/2021 data values for the row.
$ws.Range("P5").Value = 20.5
$ws.Range("Q5").Value = 20.5
$ws.Range("R5").Value = 17.9

# Match the author's final selection state (select S7:S8, active cell S7).
$ws.Range("S7:S8").Select()
